# Updated cryptos list - price (D) and volume-1h (E) refresh, plus a rank swap in rows 48/49.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    # Force the cell to remain a text value (many of the new values parse as
    # plain numbers, e.g. "243.06"), matching the source data which stores
    # prices/volumes as text, then restore the default "Normal" cell style so
    # no stray number-format is left behind.
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.174.66"
Set-TextValue $ws.Range("E2") "  +0.14%  "
Set-TextValue $ws.Range("D3") "1.829.81"
Set-TextValue $ws.Range("E3") "  -0.45%  "
Set-TextValue $ws.Range("D4") "0.9991"
Set-TextValue $ws.Range("E4") "  -0.29%  "
Set-TextValue $ws.Range("D5") "243.06"
Set-TextValue $ws.Range("E5") "  -0.03%  "
Set-TextValue $ws.Range("D6") "0.6189"
Set-TextValue $ws.Range("E6") "  +0.24%  "
Set-TextValue $ws.Range("D7") "1.001"
Set-TextValue $ws.Range("E7") "  -0.32%  "
Set-TextValue $ws.Range("D8") "0.07343"
Set-TextValue $ws.Range("E8") "  -1.77%  "
Set-TextValue $ws.Range("D9") "0.2890"
Set-TextValue $ws.Range("E9") "  -1.47%  "
Set-TextValue $ws.Range("D10") "23.21"
Set-TextValue $ws.Range("E10") "  +0.32%  "
Set-TextValue $ws.Range("D11") "0.07630"
Set-TextValue $ws.Range("E11") "  -0.94%  "
Set-TextValue $ws.Range("D12") "1.830.87"
Set-TextValue $ws.Range("E12") "  -0.04%  "
Set-TextValue $ws.Range("D13") "4.972"
Set-TextValue $ws.Range("E13") "  -0.64%  "
Set-TextValue $ws.Range("D14") "0.6695"
Set-TextValue $ws.Range("E14") "  -0.53%  "
Set-TextValue $ws.Range("E15") "  -0.56%  "
Set-TextValue $ws.Range("D16") "0.000008984"
Set-TextValue $ws.Range("E16") "  -1.80%  "
Set-TextValue $ws.Range("D17") "5.836"
Set-TextValue $ws.Range("E17") "  -1.42%  "
Set-TextValue $ws.Range("D18") "29.160.47"
Set-TextValue $ws.Range("E18") "  +0.22%  "
Set-TextValue $ws.Range("D19") "2.084.42"
Set-TextValue $ws.Range("E19") "  +0.19%  "
Set-TextValue $ws.Range("D20") "236.08"
Set-TextValue $ws.Range("E20") "  +1.12%  "
Set-TextValue $ws.Range("E21") "  -1.71%  "
Set-TextValue $ws.Range("D22") "1.001"
Set-TextValue $ws.Range("E22") "  -0.35%  "
Set-TextValue $ws.Range("D23") "7.346"
Set-TextValue $ws.Range("E23") "  +2.13%  "
Set-TextValue $ws.Range("E24") "  -0.33%  "
Set-TextValue $ws.Range("D25") "158.47"
Set-TextValue $ws.Range("E25") "  -0.57%  "
Set-TextValue $ws.Range("D26") "0.1391"
Set-TextValue $ws.Range("E26") "  -1.03%  "
Set-TextValue $ws.Range("D27") "8.517"
Set-TextValue $ws.Range("E27") "  +0.20%  "
Set-TextValue $ws.Range("D28") "17.63"
Set-TextValue $ws.Range("E28") "  -1.34%  "
Set-TextValue $ws.Range("E29") "  -0.96%  "
Set-TextValue $ws.Range("D30") "0.05855"
Set-TextValue $ws.Range("E30") "  +6.22%  "
Set-TextValue $ws.Range("D31") "1.234"
Set-TextValue $ws.Range("E31") "  +1.95%  "
Set-TextValue $ws.Range("D32") "4.080"
Set-TextValue $ws.Range("E32") "  -0.78%  "
Set-TextValue $ws.Range("D33") "4.086"
Set-TextValue $ws.Range("E33") "  -1.77%  "
Set-TextValue $ws.Range("D34") "1.864"
Set-TextValue $ws.Range("E34") "  +1.57%  "
Set-TextValue $ws.Range("E35") "  -0.52%  "
Set-TextValue $ws.Range("D36") "0.7248"
Set-TextValue $ws.Range("E36") "  -1.73%  "
Set-TextValue $ws.Range("D37") "2.609"
Set-TextValue $ws.Range("E37") "  -2.15%  "
Set-TextValue $ws.Range("D38") "2.855"
Set-TextValue $ws.Range("E38") "  +2.85%  "
Set-TextValue $ws.Range("D39") "1.228.78"
Set-TextValue $ws.Range("E39") "  +1.33%  "
Set-TextValue $ws.Range("E40") "  -1.32%  "
Set-TextValue $ws.Range("D41") "6.224"
Set-TextValue $ws.Range("E41") "  -3.58%  "
Set-TextValue $ws.Range("D42") "0.9080"
Set-TextValue $ws.Range("E42") "  +1.75%  "
Set-TextValue $ws.Range("D43") "1.001"
Set-TextValue $ws.Range("E43") "  -0.20%  "
Set-TextValue $ws.Range("D44") "102.05"
Set-TextValue $ws.Range("E44") "  +0.00%  "
Set-TextValue $ws.Range("D45") "1.987.37"
Set-TextValue $ws.Range("E45") "  +0.40%  "
Set-TextValue $ws.Range("D46") "65.64"
Set-TextValue $ws.Range("E46") "  +0.16%  "
Set-TextValue $ws.Range("D47") "0.5043"
Set-TextValue $ws.Range("E47") "  -1.12%  "
Set-TextValue $ws.Range("E50") "  -3.22%  "
Set-TextValue $ws.Range("E51") "  +2.92%  "

# Rows 48/49: TheSandbox and EnergySwap swapped ranking positions.
Set-TextValue $ws.Range("B48") "TheSandbox"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D48") "0.4045"
Set-TextValue $ws.Range("E48") "  -0.81%  "

Set-TextValue $ws.Range("B49") "EnergySwap"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "9.161"
Set-TextValue $ws.Range("E49") "  -0.10%  "
